# Append new scrape run (2026-01-18 18:25:34 JST) to the "ランサーズ" sheet.
# - refresh the timestamp on the rows that are still present
# - insert a brand-new listing ahead of the existing row 4
# - append a brand-new listing at the end (new row 7)
# - keep each URL cell's hyperlink target in sync with its new row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$stamp = "2026-01-18 18:25:34"

# --- row 2 (unchanged listing, timestamp refresh only) ---
$ws.Cells.Item(2, 1).Value = $stamp

# --- row 3 (unchanged listing, timestamp refresh only) ---
$ws.Cells.Item(3, 1).Value = $stamp

# --- row 4: NEW listing (medical-industry web app, full-stack dev) ---
$ws.Cells.Item(4, 1).Value = $stamp
$ws.Cells.Item(4, 2).Value = "【医療機関向け業務改善サービスの新規開発】WEBアプリ開発におけるフルスタック開発担当者募集"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5473940"
$ws.Cells.Item(4, 7).Value = 135
$ws.Cells.Item(4, 8).Value = "◆開発 ◇業務改善"

# --- row 5: the former row-4 listing (manufacturing DX PM), shifted down ---
$ws.Cells.Item(5, 1).Value = $stamp
$ws.Cells.Item(5, 2).Value = "製造業DXプロダクト開発のプロダクトマネージャー募集"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5468432"
$ws.Cells.Item(5, 7).Value = 75
$ws.Cells.Item(5, 8).Value = "◆開発"

# --- row 6: the former row-5 listing (urgent program fix), shifted down ---
$ws.Cells.Item(6, 1).Value = $stamp
$ws.Cells.Item(6, 2).Value = "【急募】プログラム修正依頼!スキルを活かしてみませんか?"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5473840"
$ws.Cells.Item(6, 7).Value = 13

# --- row 7: NEW listing (initial PM/PMO, partly remote) ---
$ws.Cells.Item(7, 1).Value = $stamp
$ws.Cells.Item(7, 2).Value = "初回 PM/PMO(オープン) 一部リモート"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "1,000 ~ 5,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5473958"
$ws.Cells.Item(7, 7).Value = 10

# --- rebuild the hyperlinks so every F-cell points at the right URL ---
# (Hyperlinks.Delete on any range clears the whole sheet collection in this
#  engine, so wipe once up front and re-add all six in final form.)
$ws.Range("F2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5473648")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5473858")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5473940")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5468432")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5473840")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5473958")

# Hyperlinks.Add reapplies its own font style; force the cells back onto the
# workbook's existing "Hyperlink" cell style so we don't grow the style table.
$ws.Range("F2:F7").Style = "Hyperlink"
